$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.759.31"
$ws.Range("E2").Value = "  +4.21%  "
$ws.Range("D3").Value = "3.072.01"
$ws.Range("E3").Value = "  +2.60%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.59"
$ws.Range("E5").Value = "  +2.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.54"
$ws.Range("E6").Value = "  +2.74%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "3.061.66"
$ws.Range("E8").Value = "  +2.65%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("E10").Value = "  +5.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.67"
$ws.Range("E11").Value = "  +10.76%  "
$ws.Range("E12").Value = "  +1.94%  "
$ws.Range("E13").Value = "  +4.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.10"
$ws.Range("E14").Value = "  +4.37%  "
$ws.Range("E15").Value = "  +0.48%  "
$ws.Range("D16").Value = "3.580.23"
$ws.Range("E16").Value = "  +2.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.23"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").Value = "3.070.84"
$ws.Range("E18").Value = "  +2.68%  "
$ws.Range("D19").Value = "61.696.27"
$ws.Range("E19").Value = "  +4.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "446.15"
$ws.Range("E20").Value = "  +3.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.97"
$ws.Range("E21").Value = "  +2.23%  "
$ws.Range("E22").Value = "  +1.89%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.42"
$ws.Range("E23").Value = "  +4.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.71"
$ws.Range("E24").Value = "  +2.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.61"
$ws.Range("E25").Value = "  +0.89%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  +5.53%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("E29").Value = "  +4.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.18"
$ws.Range("E30").Value = "  +5.68%  "
$ws.Range("E31").Value = "  +10.51%  "
$ws.Range("E32").Value = "  +14.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.77"
$ws.Range("E33").Value = "  +4.09%  "
$ws.Range("E34").Value = "  +4.11%  "
$ws.Range("E35").Value = "  +3.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.03"
$ws.Range("E36").Value = "  +1.87%  "
$ws.Range("E37").Value = "  +5.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "50.14"
$ws.Range("E38").Value = "  +2.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.94"
$ws.Range("E39").Value = "  +9.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.80"
$ws.Range("E40").Value = "  +1.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "419.72"
$ws.Range("E41").Value = "  +4.72%  "
$ws.Range("D42").Value = "2.963.26"
$ws.Range("E42").Value = "  +7.57%  "
$ws.Range("E43").Value = "  +4.99%  "
$ws.Range("E44").Value = "  +10.06%  "
$ws.Range("E45").Value = "  +0.43%  "
$ws.Range("E46").Value = "  +5.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.46"
$ws.Range("E48").Value = "  +2.60%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.86"
$ws.Range("E49").Value = "  +0.39%  "
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.30"
$ws.Range("E51").Value = "  +3.98%  "
